# Add three new columns (color_1, color_2, agency_level) to the conditions
# table, fix a bad "group" value in row 16, and fill in all the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header row -------------------------------------------------------
$ws.Range("F1").Value = "color_1"
$ws.Range("G1").Value = "color_2"
$ws.Range("H1").Value = "agency_level"

# Give the new header cells the same look as the rest of the header row.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial($xlPasteFormats)

# --- Fix existing data: row 16's "group" should be 2, not 3 -----------
$ws.Range("E16").Value = 2

# --- New body data (rows 2-17) -----------------------------------------
$rows = @(
    @{ F = "red";    G = "blue";   H = "high" },
    @{ F = "red";    G = "green";  H = "high" },
    @{ F = "yellow"; G = "grey";   H = "high" },
    @{ F = "yellow"; G = "purple"; H = "low"  },
    @{ F = "red";    G = "green";  H = "low"  },
    @{ F = "red";    G = "blue";   H = "low"  },
    @{ F = "blue";   G = "white";  H = "high" },
    @{ F = "green";  G = "blue";   H = "low"  },
    @{ F = "red";    G = "blue";   H = "high" },
    @{ F = "white";  G = "white";  H = "high" },
    @{ F = "blue";   G = "blue";   H = "low"  },
    @{ F = "yellow"; G = "white";  H = "high" },
    @{ F = "red";    G = "green";  H = "low"  },
    @{ F = "red";    G = "green";  H = "high" },
    @{ F = "red";    G = "green";  H = "low"  },
    @{ F = "red";    G = "blue";   H = "low"  }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H

    # Match the formatting already used by the other text cells (B:D) in
    # this row, rather than leaving the default/plain style behind.
    $ws.Range("B$r").Copy()
    $ws.Range("F${r}:H$r").PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false
